# Generate Report for Handback
# Adds a new tracked file (c72a3c4e-96b5-44ec-92b2-72eec808df73.md) as row 4
# on the "Overview", "zh-cn" and "de-de" sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$dateCreated   = "2016-08-15 16:42:50"
$zhGenDate     = "2016-08-15 16:42:46"
$zhHandback    = "2016-08-15 16:43:06"
$deGenDate     = "2016-08-15 16:42:50"
$deHandback    = "2016-08-15 16:43:14"
$dtFormat      = "yyyy-mm-dd HH:mm:ss"

$fileName      = "c72a3c4e-96b5-44ec-92b2-72eec808df73.md"
$pathAndName   = "e2e\c72a3c4e-96b5-44ec-92b2-72eec808df73.md"
$status        = "Handed back: in sync with en-US"
$zhXlf         = "c72a3c4e-96b5-44ec-92b2-72eec808df73.021def19692dd357980ccdc6f4dd6fa81cb66694.zh-cn.xlf"
$deXlf         = "c72a3c4e-96b5-44ec-92b2-72eec808df73.021def19692dd357980ccdc6f4dd6fa81cb66694.de-de.xlf"

$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/021def19692dd357980ccdc6f4dd6fa81cb66694/e2e/c72a3c4e-96b5-44ec-92b2-72eec808df73.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/021def19692dd357980ccdc6f4dd6fa81cb66694/e2e/c72a3c4e-96b5-44ec-92b2-72eec808df73.md"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/021def19692dd357980ccdc6f4dd6fa81cb66694/e2e/c72a3c4e-96b5-44ec-92b2-72eec808df73.md"

# ===================== Sheet "Overview" =====================
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A4").Value = $fileName
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = $status
$wsOv.Range("F4").Value = $status
$wsOv.Range("G4").Value = $dateCreated
$wsOv.Range("G4").NumberFormat = $dtFormat

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), $mdUrl, "", "", $pathAndName)

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G4"))

# ===================== Sheet "zh-cn" =====================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhGenDate
$wsZh.Range("H4").NumberFormat = $dtFormat
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHandback
$wsZh.Range("K4").NumberFormat = $dtFormat
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhUrl, "", "", $fileName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ===================== Sheet "de-de" =====================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deGenDate
$wsDe.Range("H4").NumberFormat = $dtFormat
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHandback
$wsDe.Range("K4").NumberFormat = $dtFormat
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deUrl, "", "", $fileName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
